$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the number format used on the existing D9 date cell so the
# duplicated record (moved to row 10) keeps the same formatting.
$dateFormat = $ws.Range("D9").NumberFormat

# Copy the current (old) row 9 values down to row 10 first, since the new
# row 9 data pushes the old row 9 record to row 10.
for ($col = 1; $col -le 20; $col++) {
    $ws.Cells.Item(10, $col).Value = $ws.Cells.Item(9, $col).Value2
}
$ws.Range("D10").NumberFormat = $dateFormat

# Now update row 9 with the new record values.
$ws.Range("D9").Value = 44511
$ws.Range("M9").Value = 120
$ws.Range("N9").Value = 28000
$ws.Range("O9").Value = 28000
$ws.Range("P9").Value = 28000
$ws.Range("S9").Value = 2800
